# Front End Tricks.xlsx - "HTML-CSS" sheet
# Adds 6 new rows (50-55) of HTML/CSS tricks, shrinks every data row's
# height from 120.75 to 93, and moves the active selection/top-left view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HTML-CSS")

$xlPasteFormats = -4122

function Add-TrickRow {
    param($Row, $FormatSourceRow, $BText, $CText, $DText, $DUrl)

    # B (term) and D (source link) cells don't exist yet on these rows, so
    # clone formatting from a row that already has the desired style, then
    # stamp in the literal value.
    $ws.Range("B$FormatSourceRow").Copy() | Out-Null
    $ws.Range("B$Row").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("B$Row").Value = $BText

    # C (example) cell already exists as a blank, correctly-styled
    # placeholder cell on every one of these rows - just set its value.
    $ws.Range("C$Row").Value = $CText

    # D (source link): set the text + hyperlink first, then re-apply the
    # cloned formatting last, since Hyperlinks.Add() stamps its own
    # built-in "Hyperlink" look on the cell when it runs.
    $ws.Range("D$Row").Value = $DText
    $ws.Hyperlinks.Add($ws.Range("D$Row"), $DUrl) | Out-Null
    $ws.Range("D$FormatSourceRow").Copy() | Out-Null
    $ws.Range("D$Row").PasteSpecial($xlPasteFormats) | Out-Null

    $excel.CutCopyMode = 0
}

Add-TrickRow 50 30 "<body>" `
    "The <body> tag defines the document's body. The <body> element contains all the contents of an HTML document, such as text, hyperlinks, images, tables, lists, etc." `
    "https://www.w3schools.com/tags/tag_body.asp " `
    "https://www.w3schools.com/tags/tag_body.asp"

Add-TrickRow 51 30 "CSS Inheritance" `
    "Inheritance in CSS is the mechanism through which certain properties are passed on from a parent element down to its children." `
    "https://www.w3.org/wiki/Inheritance_and_cascade " `
    "https://www.w3.org/wiki/Inheritance_and_cascade"

$overrideSpacer = "".PadRight(83)
$overrideStyle = "CSS Override Style:" + $overrideSpacer + ".pink-text {`n    color: pink !important;`n  }"
Add-TrickRow 52 49 $overrideStyle `
    "Class overrides parent style; Second class overrides first class; ID overrides class; Inline overrides ID; !important overrides all" `
    "https://zenorocha.com/css-important/ " `
    "https://zenorocha.com/css-important/"

Add-TrickRow 53 49 "Colors: Hex Code" `
    "Examples: #ff0000, #ccffff, #660066" `
    "https://www.w3schools.com/colors/colors_picker.asp" `
    "https://www.w3schools.com/colors/colors_picker.asp"

Add-TrickRow 54 30 "Short Hex Codes" `
    "Example: #ff6600 - > #f60" `
    "http://www.websiteoptimization.com/speed/tweak/hex/ " `
    "http://www.websiteoptimization.com/speed/tweak/hex/"

Add-TrickRow 55 30 "Colors: RGB" `
    "Examples: rgb (255, 0, 0), rgb (218, 112, 214), rgb (160, 82, 45)" `
    "https://www.w3schools.com/colors/colors_rgb.asp" `
    "https://www.w3schools.com/colors/colors_rgb.asp"

# All data rows (2-131) shrink from 120.75pt to 93pt.
$ws.Rows("2:131").RowHeight = 93

# Move the view: top-left cell scrolls to row 54, and the new selection
# lands on D56 (the first still-empty row right after the new content).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D56").Select() | Out-Null
